# Apply the "large set of changes" commit to the FAQ workbook.
#
# Net effect on xl/sharedStrings.xml: the four old pricing / contact-us
# strings (Finnish + English, in both the "tilaus" FAQ and its English
# translation) are replaced by updated copy (17e/h instead of 20e/h for
# design, 30/35e/h instead of 40e/h for coding, and a reworded "tight
# budget" condition instead of "clear needs/tarpeesi ovat selkeät").
# Updating just these four cells is enough - the engine prunes now-unused
# shared strings and appends the new ones on save, which reproduces the
# reindex seen across every other worksheet in the diff.
#
# Also: the active sheet/tab moves from "faq-other" to "en-tilaus", with
# new cell selections on a couple of sheets.

$wb = $excel.ActiveWorkbook

$faqTilaus = $wb.Worksheets.Item("faq-tilaus")
$enTilaus  = $wb.Worksheets.Item("en-tilaus")
$faqOther  = $wb.Worksheets.Item("faq-other")

# --- Updated copy -----------------------------------------------------

$faqTilausContact = "Ota yhteyttä ja kerro minkälainen projekti on kysessä. Kartoitan projektin työnmäärän ja teen tarjouksen. Tarjoukseen on määritelty aika arvio työvaiheista ja projektin lopullisesta deadlinesta, jos budjettisi on tiukka voimme myös sopia kiinteästä hinnasta. Tarjouksen pyytäminen ja siinä tapahtuva konsultointi on täysin ilmaista."

$faqTilausPricing = "laskutan 17e/h suunnittelusta ja noin 35e/h koodaamisesta riippuen mitä tekniikoita käytän. Tällöin yksinkertaiset verkkosivut maksavat noin 300e - 750e ja verkkopalvelujen tuottaminen noin 1200e - 3000e. Jos sinulla on tiukka budjetti voimme myös sopia kiinteästä hinnasta ja tehdä sivuston hostaaminen 0e/kk."

$enTilausPricing = "I charge 17e/h for the design and about 30e/h coding, depending on what technologies I use. In this case, the simple webpages cost about 300e - 750e and the complete web service is about 1200e - 3000e. If you have a tight budget we can also agree on a fixed price and make the sites hosting cost 0e/month."

$enTilausContact = "Contact me and tell what kind of project is in question. I'll map the workload of the project and will make you an offer.`nIf you budget is tight we can also agree on a fixed price. Requesting the offer is completely free."

# --- Cell edits ---------------------------------------------------------
# Order matters: new shared strings are appended in the order they are
# first written, and the diff expects laskutan17/Icharge17/OtaYhteytta/
# ContactMe to land consecutively at the tail of the table.

$faqTilaus.Range("B3").Value = $faqTilausPricing
$enTilaus.Range("B3").Value = $enTilausPricing
$faqTilaus.Range("B2").Value = $faqTilausContact
$enTilaus.Range("B2").Value = $enTilausContact

# Writing the longer English contact copy auto-grows row 2's custom
# height; restore its original explicit height (15pt) so only the text
# itself changes.
$enTilaus.Rows.Item(2).RowHeight = 15

# --- Selections -----------------------------------------------------

$faqTilaus.Range("B4").Select()
$faqOther.Range("B2").Select()

# en-tilaus becomes the active sheet/tab, selection B2 - select this last
# so it ends up as the active tab on save.
$enTilaus.Range("B2").Select()
